$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "AddCounterparty" sheet right before "FlagReason".
# ---------------------------------------------------------------------------
$flagReason = $wb.Worksheets.Item("FlagReason")
$addCounterparty = $wb.Worksheets.Add($flagReason)
$addCounterparty.Name = "AddCounterparty"

$addCounterparty.Range("A1").Value = "Company Name"
$addCounterparty.Range("B1").Value = "Type"
$addCounterparty.Range("A1:B1").Font.Bold = $true

$addCounterparty.Range("A2").Value = "Zillow Home Loans, LLC"
$addCounterparty.Range("B2").Value = "Financial"

$addCounterparty.Columns.Item(1).ColumnWidth = 20.4
$addCounterparty.Columns.Item(2).ColumnWidth = 21

$addCounterparty.Range("D10").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2) Insert the new "Bid" sheet right after "FlagReason".
# ---------------------------------------------------------------------------
$flagReason = $wb.Worksheets.Item("FlagReason")
$bid = $wb.Worksheets.Add($null, $flagReason)
$bid.Name = "Bid"

$bid.Range("A1").Value = "Round Name"
$bid.Range("B1").Value = "Amount"
$bid.Range("A1:B1").Font.Bold = $true
$bid.Range("A1:B1").HorizontalAlignment = -4108
$bid.Range("A1:B1").VerticalAlignment = -4108

$bid.Range("A2").Value = "Closing"

$bid.Range("B2").HorizontalAlignment = -4131
$bid.Range("B2").VerticalAlignment = -4160
$bid.Range("B2").NumberFormat = "@"
$bid.Range("B2").Value = "10"

$bid.Columns.Item(1).ColumnWidth = 11.3

$bid.Activate() | Out-Null
$bid.Range("E12").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3) Misc view-state tweaks carried by the diff.
# ---------------------------------------------------------------------------
# AddOpportunity: scroll so column F is the leftmost visible column and
# selection moves to AB2.
$addOpportunity = $wb.Worksheets.Item("AddOpportunity")
$addOpportunity.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 6
$addOpportunity.Range("AB2").Select() | Out-Null

# Re-activate Bid last so it ends up the selected tab, matching the new
# activeTab or the workbook.
$bid.Activate() | Out-Null
